$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unify the "strain" column (F) formatting for rows 3-27 to match row 2's
# "KN99alpha" (no space), so all metadata rows use the same format.
$ws.Range("F3:F27").Value = "KN99alpha"

# Match the saved selection state recorded in the workbook.
$ws.Range("F3:F27").Select()
